# Update benchmark data in the "chart" worksheet (the active sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("chart")

$data = @(
    @(3984, 3962, 312, 1588, 1289, 324),
    @(4878, 4863, 303, 1522, 1279, 312),
    @(4228, 4018, 570, 2388, 1513, 913),
    @(4078, 3989, 393, 1626, 1405, 257),
    @(5044, 4891, 401, 1557, 1343, 259),
    @(4169, 3908, 563, 2253, 1561, 706),
    @(4038, 3869, 423, 1382, 1105, 316),
    @(4096, 3968, 417, 1431, 1144, 311),
    @(4843, 4025, 1085, 2610, 1605, 1066),
    @(5281, 4871, 707, 1893, 1331, 577),
    @(5107, 4858, 508, 1726, 1324, 417),
    @(4197, 3871, 605, 1517, 1102, 423),
    @(5158, 4884, 515, 1890, 1480, 449),
    @(5289, 4931, 664, 1965, 1521, 521)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

$excel.Calculate()

# Force chart caches to refresh by touching each chart object on the sheet
foreach ($co in $ws.ChartObjects()) {
    $co.Chart.Refresh()
}

$wb.Save()
